# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1312962.6
$ws.Range("I12").Value = 1429100.1
$ws.Range("J12").Value = 500000
$ws.Range("K12").Value = 1429100.1
$ws.Range("L12").Value = 500000
$ws.Range("M12").Value = -1428930.1
$ws.Range("N12").Value = -500340
$ws.Range("H33").Value = 62723.312
$ws.Range("I33").Value = 77183.84
$ws.Range("J33").Value = 61
$ws.Range("K33").Value = 77183.84
$ws.Range("L33").Value = 61
$ws.Range("M33").Value = -76954.84
$ws.Range("N33").Value = -519
$ws.Range("H74").Value = 3057507
$ws.Range("I74").Value = 3996663.2
$ws.Range("K74").Value = 3996663.2
$ws.Range("M74").Value = -3995727.2
$ws.Range("H77").Value = 3057507
$ws.Range("I77").Value = 3996663.2
$ws.Range("K77").Value = 19983316
$ws.Range("M77").Value = -19978636
$ws.Range("H92").Value = 1755.75
$ws.Range("I92").Value = 1823.2307
$ws.Range("K92").Value = 1823.2307
$ws.Range("M92").Value = -575.2307000000001
$ws.Range("H116").Value = 3735.8462
$ws.Range("I116").Value = 4026
$ws.Range("J116").Value = 2768.6667
$ws.Range("K116").Value = 4026
$ws.Range("L116").Value = 2768.6667
$ws.Range("M116").Value = -584
$ws.Range("N116").Value = -9652.6667
$ws.Range("H138").Value = 2594.71
$ws.Range("J138").Value = 3298.822
$ws.Range("L138").Value = 9896.466
$ws.Range("N138").Value = -20176.466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4081.012
$ws.Range("I32").Value = 3667.2817
$ws.Range("J32").Value = 6340.615
$ws.Range("K32").Value = 3667.2817
$ws.Range("L32").Value = 6340.615
$ws.Range("M32").Value = -3380.2817
$ws.Range("N32").Value = -6914.615
$ws.Range("H122").Value = 1966.3334
$ws.Range("I122").Value = 1966.3334
$ws.Range("K122").Value = 5899.0002
$ws.Range("M122").Value = -3449.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2341.6843
$ws.Range("I86").Value = 2072.8
$ws.Range("J86").Value = 3350
$ws.Range("K86").Value = 2072.8
$ws.Range("L86").Value = 3350
$ws.Range("M86").Value = -949.8000000000002
$ws.Range("N86").Value = -5596
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812
$ws.Range("H89").Value = 2341.6843
$ws.Range("I89").Value = 2072.8
$ws.Range("J89").Value = 3350
$ws.Range("K89").Value = 10364
$ws.Range("L89").Value = 16750
$ws.Range("M89").Value = -4748
$ws.Range("N89").Value = -27982
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1846.7317
$ws.Range("I31").Value = 1846.7317
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1846.7317
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -1551.7317
$ws.Range("H34").Value = 1846.7317
$ws.Range("I34").Value = 1846.7317
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1846.7317
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -1644.7317
$ws.Range("H58").Value = 3375.3333
$ws.Range("I58").Value = 562.5517
$ws.Range("K58").Value = 562.5517
$ws.Range("M58").Value = -359.5517
$ws.Range("H94").Value = 5093.1113
$ws.Range("J94").Value = 4975.143
$ws.Range("L94").Value = 4975.143
$ws.Range("N94").Value = -5877.143
$ws.Range("H134").Value = 4633.4165
$ws.Range("I134").Value = 3475.125
$ws.Range("J134").Value = 6950
$ws.Range("K134").Value = 10425.375
$ws.Range("L134").Value = 20850
$ws.Range("M134").Value = -7890.375
$ws.Range("N134").Value = -25920
$ws.Range("H136").Value = 3375.3333
$ws.Range("I136").Value = 562.5517
$ws.Range("K136").Value = 1687.6551
$ws.Range("M136").Value = 862.3449000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1348.8
$ws.Range("I5").Value = 765.3333
$ws.Range("K5").Value = 2295.9999
$ws.Range("M5").Value = -2183.9999
$ws.Range("H131").Value = 32461.637
$ws.Range("J131").Value = 2542.037
$ws.Range("L131").Value = 7626.110999999999
$ws.Range("N131").Value = -17706.111
$ws.Range("H135").Value = 1348.8
$ws.Range("I135").Value = 765.3333
$ws.Range("K135").Value = 6887.9997
$ws.Range("M135").Value = -4352.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4184.8887
$ws.Range("I70").Value = 3900.3333
$ws.Range("J70").Value = 5323.1113
$ws.Range("K70").Value = 3900.3333
$ws.Range("L70").Value = 5323.1113
$ws.Range("M70").Value = -3630.3333
$ws.Range("N70").Value = -5863.1113
$ws.Range("H73").Value = 4184.8887
$ws.Range("I73").Value = 3900.3333
$ws.Range("J73").Value = 5323.1113
$ws.Range("K73").Value = 3900.3333
$ws.Range("L73").Value = 5323.1113
$ws.Range("M73").Value = -2964.3333
$ws.Range("N73").Value = -7195.1113
$ws.Range("H126").Value = 2520.8
$ws.Range("I126").Value = 2545.2727
$ws.Range("J126").Value = 2453.5
$ws.Range("K126").Value = 7635.8181
$ws.Range("L126").Value = 7360.5
$ws.Range("M126").Value = -5165.8181
$ws.Range("N126").Value = -12300.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 469.23077
$ws.Range("I46").Value = 450
$ws.Range("J46").Value = 477.77777
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 477.77777
$ws.Range("M46").Value = -262
$ws.Range("N46").Value = -853.7777699999999
$ws.Range("H68").Value = 2816.8333
$ws.Range("I68").Value = 2880.4
$ws.Range("J68").Value = 2771.4285
$ws.Range("K68").Value = 2880.4
$ws.Range("L68").Value = 2771.4285
$ws.Range("M68").Value = -2131.4
$ws.Range("N68").Value = -4269.4285
$ws.Range("H71").Value = 2816.8333
$ws.Range("I71").Value = 2880.4
$ws.Range("J71").Value = 2771.4285
$ws.Range("K71").Value = 14402
$ws.Range("L71").Value = 13857.1425
$ws.Range("M71").Value = -10658
$ws.Range("N71").Value = -21345.1425
$ws.Range("H82").Value = 2543.0715
$ws.Range("I82").Value = 2362.5
$ws.Range("J82").Value = 2783.8333
$ws.Range("K82").Value = 2362.5
$ws.Range("L82").Value = 2783.8333
$ws.Range("M82").Value = -2001.5
$ws.Range("N82").Value = -3505.8333
$ws.Range("H85").Value = 2543.0715
$ws.Range("I85").Value = 2362.5
$ws.Range("J85").Value = 2783.8333
$ws.Range("K85").Value = 2362.5
$ws.Range("L85").Value = 2783.8333
$ws.Range("M85").Value = -1114.5
$ws.Range("N85").Value = -5279.8333
$ws.Range("H136").Value = 2534.3333
$ws.Range("I136").Value = 1405.2858
$ws.Range("J136").Value = 8857
$ws.Range("K136").Value = 4215.857400000001
$ws.Range("L136").Value = 26571
$ws.Range("M136").Value = -1665.857400000001
$ws.Range("N136").Value = -31671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 227.60869
$ws.Range("I113").Value = 216.77777
$ws.Range("J113").Value = 266.6
$ws.Range("K113").Value = 650.33331
$ws.Range("L113").Value = 799.8000000000001
$ws.Range("M113").Value = 1519.66669
$ws.Range("N113").Value = -5139.8
$ws.Range("H136").Value = 2249.8518
$ws.Range("I136").Value = 2122.4
$ws.Range("J136").Value = 2409.1667
$ws.Range("K136").Value = 6367.200000000001
$ws.Range("L136").Value = 7227.500100000001
$ws.Range("M136").Value = -3817.200000000001
$ws.Range("N136").Value = -12327.5001
